# This script applies updated Betfair Back/Lay odds values to Sheet1,
# cell by cell, matching the authoritative diff of the edited workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 5  # was 4.8
$ws.Range("G2").Value = 5.3  # was 5.1
$ws.Range("H2").Value = 1.85  # was 1.86
$ws.Range("I2").Value = 1.88  # was 1.89
$ws.Range("R2").Value = 1.3  # was 1.29
$ws.Range("T2").Value = 2.04  # was 2.02
$ws.Range("U2").Value = 1.9  # was 1.89
$ws.Range("W2").Value = 1.23  # was 1.24
$ws.Range("Y2").Value = 7.6  # was 7.8
$ws.Range("AB2").Value = 15  # was 14.5
$ws.Range("AG2").Value = 21  # was 19.5
$ws.Range("AI2").Value = 44  # was 48
$ws.Range("AJ2").Value = 1000  # was 130
$ws.Range("AK2").Value = 80  # was 75
$ws.Range("AM2").Value = 150  # was 140
$ws.Range("AN2").Value = 120  # was 110
$ws.Range("AO2").Value = 15  # was 15.5

# Row 3
$ws.Range("F3").Value = 1.05  # was 1.04
$ws.Range("H3").Value = 1.05  # was 1.04
$ws.Range("J3").Value = 3.65  # was 1.04
$ws.Range("K3").Value = 950  # was 1000
$ws.Range("M3").Value = 1.02  # was 1.01
$ws.Range("N3").Value = 1.31  # was 1.3
$ws.Range("T3").Value = 1.05  # was 1.04
$ws.Range("V3").Value = 1.02  # was 1.01
$ws.Range("W3").Value = 1.02  # was 1.01

# Row 4
$ws.Range("G4").Value = 600  # was 610
$ws.Range("J4").Value = 3.6  # was 1.02
$ws.Range("M4").Value = 1.02  # was 1.01
$ws.Range("N4").Value = 1.26  # was 1.25
$ws.Range("P4").Value = 1.26  # was 1.24
$ws.Range("R4").Value = 1.21  # was 1.18
$ws.Range("V4").Value = 1.02  # was 1.01
$ws.Range("W4").Value = 1.02  # was 1.01

# Row 5
$ws.Range("G5").Value = 1.53  # was 1.55
$ws.Range("I5").Value = 9.2  # was 8.6
$ws.Range("L5").Value = 1.13  # was 1.01
$ws.Range("O5").Value = 1.27  # was 1.28
$ws.Range("P5").Value = 2.1  # was 2.06
$ws.Range("S5").Value = 3  # was 2.96
$ws.Range("T5").Value = 1.97  # was 1.94
$ws.Range("W5").Value = 2.88  # was 2.8
$ws.Range("AB5").Value = 8.6  # was 8.8
$ws.Range("AF5").Value = 9.2  # was 10.5
$ws.Range("AJ5").Value = 1000  # was 14.5
$ws.Range("AN5").Value = 8  # was 9

# Row 6
$ws.Range("J6").Value = 3.75  # was 3.25
$ws.Range("K6").Value = 950  # was 25
$ws.Range("M6").Value = 1.02  # was 1.01
$ws.Range("N6").Value = 3.15  # was 3.05
$ws.Range("P6").Value = 3.15  # was 3.05
$ws.Range("Q6").Value = 1.34  # was 1.37
$ws.Range("R6").Value = 2.16  # was 2.04
$ws.Range("S6").Value = 1.81  # was 1.83
$ws.Range("T6").Value = 1.31  # was 1.05
$ws.Range("W6").Value = 1.33  # was 1.28

# Row 7
$ws.Range("K7").Value = 3.8  # was 3.85
$ws.Range("M7").Value = 1.07  # was 1.06
$ws.Range("P7").Value = 1.98  # was 1.96
$ws.Range("Q7").Value = 1.96  # was 1.97
$ws.Range("R7").Value = 1.39  # was 1.38
$ws.Range("S7").Value = 3.4  # was 3.45

# Row 8
$ws.Range("F8").Value = 1.35  # was 1.34
$ws.Range("G8").Value = 1.36  # was 1.35
$ws.Range("J8").Value = 5.5  # was 5.6
$ws.Range("K8").Value = 5.6  # was 5.7
$ws.Range("N8").Value = 4  # was 3.9
$ws.Range("O8").Value = 1.31  # was 1.32
$ws.Range("Q8").Value = 1.96  # was 1.97
$ws.Range("R8").Value = 1.38  # was 1.37
$ws.Range("S8").Value = 3.45  # was 3.5
$ws.Range("X8").Value = 16.5  # was 16
$ws.Range("Y8").Value = 32  # was 34
$ws.Range("AA8").Value = 600  # was 590
$ws.Range("AI8").Value = 250  # was 230
$ws.Range("AN8").Value = 6.6  # was 6.8
$ws.Range("AO8").Value = 490  # was 410

# Row 9
$ws.Range("F9").Value = 2.7  # was 2.68
$ws.Range("I9").Value = 2.82  # was 2.8
$ws.Range("K9").Value = 3.75  # was 3.85
$ws.Range("N9").Value = 3.9  # was 3.85
$ws.Range("V9").Value = 1.54  # was 1.55
$ws.Range("W9").Value = 1.5  # was 1.49
$ws.Range("Y9").Value = 13  # was 13.5
$ws.Range("AB9").Value = 13  # was 13.5
$ws.Range("AG9").Value = 13.5  # was 15.5
$ws.Range("AH9").Value = 17  # was 18.5
$ws.Range("AK9").Value = 980  # was 36
$ws.Range("AM9").Value = 95  # was 100
$ws.Range("AO9").Value = 26  # was 25

# Row 10
$ws.Range("F10").Value = 1.68  # was 1.65
$ws.Range("H10").Value = 5.3  # was 5
$ws.Range("I10").Value = 6.4  # was 6.2
$ws.Range("K10").Value = 4.3  # was 4.4
$ws.Range("N10").Value = 4.6  # was 4.2
$ws.Range("Q10").Value = 1.73  # was 1.74
$ws.Range("R10").Value = 1.46  # was 1.44
$ws.Range("S10").Value = 2.78  # was 2.84
$ws.Range("T10").Value = 1.71  # was 1.77
$ws.Range("U10").Value = 2.2  # was 2.06
$ws.Range("X10").Value = 27  # was 21
$ws.Range("Y10").Value = 24  # was 21
$ws.Range("Z10").Value = 55  # was 1000
$ws.Range("AB10").Value = 11  # was 12.5
$ws.Range("AC10").Value = 9  # was 11.5
$ws.Range("AD10").Value = 24  # was 22
$ws.Range("AE10").Value = 90  # was 1000
$ws.Range("AF10").Value = 11.5  # was 13.5
$ws.Range("AG10").Value = 11  # was 13.5
$ws.Range("AH10").Value = 18  # was 20
$ws.Range("AI10").Value = 75  # was 1000
$ws.Range("AK10").Value = 18  # was 18.5
$ws.Range("AO10").Value = 85  # was 1000

# Row 11
$ws.Range("F11").Value = 1.49  # was 1.48
$ws.Range("G11").Value = 1.55  # was 1.54
$ws.Range("H11").Value = 7  # was 6.8
$ws.Range("J11").Value = 4.6  # was 4.7
$ws.Range("K11").Value = 5.1  # was 5.2
$ws.Range("T11").Value = 1.76  # was 1.72
$ws.Range("U11").Value = 2.18  # was 2.14
$ws.Range("W11").Value = 2.8  # was 2.84
$ws.Range("X11").Value = 970  # was 28
$ws.Range("Y11").Value = 970  # was 32
$ws.Range("Z11").Value = 65  # was 75
$ws.Range("AC11").Value = 12  # was 13.5
$ws.Range("AF11").Value = 11.5  # was 13
$ws.Range("AI11").Value = 80  # was 95
$ws.Range("AK11").Value = 16  # was 18
$ws.Range("AO11").Value = 90  # was 100

# Row 12
$ws.Range("F12").Value = 2.14  # was 2.2
$ws.Range("G12").Value = 2.24  # was 2.22
$ws.Range("I12").Value = 3.65  # was 3.6
$ws.Range("J12").Value = 3.7  # was 3.75
$ws.Range("K12").Value = 3.95  # was 3.9
$ws.Range("L12").Value = 1.26  # was 1.27
$ws.Range("T12").Value = 1.6  # was 1.56
$ws.Range("V12").Value = 1.37  # was 1.38
$ws.Range("W12").Value = 1.8  # was 1.81

# Row 13
$ws.Range("F13").Value = 1.94  # was 1.91
$ws.Range("G13").Value = 1.96  # was 1.93
$ws.Range("I13").Value = 4.2  # was 4.3
$ws.Range("N13").Value = 5.3  # was 5.4
$ws.Range("P13").Value = 2.46  # was 2.44
$ws.Range("Q13").Value = 1.66  # was 1.65
$ws.Range("S13").Value = 2.52  # was 2.54
$ws.Range("T13").Value = 1.62  # was 1.63
$ws.Range("U13").Value = 2.52  # was 2.5
$ws.Range("V13").Value = 1.31  # was 1.3
$ws.Range("W13").Value = 2.04  # was 2.06
$ws.Range("AA13").Value = 80  # was 85
$ws.Range("AB13").Value = 12.5  # was 13
$ws.Range("AJ13").Value = 23  # was 22
$ws.Range("AN13").Value = 8.8  # was 8.6

# Row 14
$ws.Range("R14").Value = 2.2  # was 2.22
$ws.Range("T14").Value = 1.7  # was 1.71
$ws.Range("U14").Value = 2.36  # was 2.34

# Row 15
$ws.Range("F15").Value = 1.82  # was 1.81
$ws.Range("O15").Value = 1.21  # was 1.2
$ws.Range("P15").Value = 2.5  # was 2.52
$ws.Range("Q15").Value = 1.65  # was 1.64
$ws.Range("R15").Value = 1.6  # was 1.61
$ws.Range("S15").Value = 2.6  # was 2.54
$ws.Range("U15").Value = 2.48  # was 2.52
$ws.Range("V15").Value = 1.26  # was 1.25
$ws.Range("Y15").Value = 23  # was 24
$ws.Range("AA15").Value = 100  # was 95
$ws.Range("AC15").Value = 9.4  # was 9.6
$ws.Range("AH15").Value = 16.5  # was 16

# Row 16
$ws.Range("F16").Value = 1.17  # was 1.16
$ws.Range("G16").Value = 1.18  # was 1.17
$ws.Range("J16").Value = 9  # was 9.6
$ws.Range("K16").Value = 9.4  # was 10
$ws.Range("N16").Value = 5.8  # was 5.7
$ws.Range("P16").Value = 2.62  # was 2.64
$ws.Range("Q16").Value = 1.6  # was 1.59
$ws.Range("R16").Value = 1.63  # was 1.64
$ws.Range("S16").Value = 2.54  # was 2.52
$ws.Range("T16").Value = 2.82  # was 2.84
$ws.Range("U16").Value = 1.52  # was 1.51
$ws.Range("W16").Value = 6.6  # was 7
$ws.Range("Z16").Value = 340  # was 360
$ws.Range("AH16").Value = 70  # was 65
$ws.Range("AI16").Value = 540  # was 520
$ws.Range("AJ16").Value = 7.4  # was 7.2
$ws.Range("AN16").Value = 3.6  # was 3.55

# Row 17
$ws.Range("H17").Value = 9  # was 9.2
$ws.Range("J17").Value = 4.5  # was 4.6
$ws.Range("K17").Value = 4.7  # was 4.8
$ws.Range("N17").Value = 4  # was 3.95
$ws.Range("U17").Value = 1.83  # was 1.8
$ws.Range("W17").Value = 3.05  # was 3.1
$ws.Range("Y17").Value = 28  # was 27
$ws.Range("Z17").Value = 80  # was 85
$ws.Range("AA17").Value = 380  # was 390
$ws.Range("AC17").Value = 10.5  # was 11
$ws.Range("AI17").Value = 140  # was 150
$ws.Range("AL17").Value = 50  # was 48
$ws.Range("AN17").Value = 8  # was 7.6

# Row 18
$ws.Range("G18").Value = 7.8  # was 7.6
$ws.Range("J18").Value = 4.9  # was 5
$ws.Range("L18").Value = 1.31  # was 1.3
$ws.Range("N18").Value = 4.9  # was 5.1
$ws.Range("O18").Value = 1.22  # was 1.21
$ws.Range("P18").Value = 2.36  # was 2.46
$ws.Range("Q18").Value = 1.65  # was 1.62
$ws.Range("R18").Value = 1.53  # was 1.56
$ws.Range("S18").Value = 2.68  # was 2.56
$ws.Range("T18").Value = 1.84  # was 1.78
$ws.Range("U18").Value = 2.12  # was 2.14
$ws.Range("V18").Value = 2.96  # was 3
$ws.Range("W18").Value = 1.14  # was 1.15
$ws.Range("X18").Value = 22  # was 23
$ws.Range("Y18").Value = 9.8  # was 10.5
$ws.Range("Z18").Value = 10.5  # was 11
$ws.Range("AB18").Value = 32  # was 980
$ws.Range("AG18").Value = 28  # was 27
$ws.Range("AH18").Value = 22  # was 21
$ws.Range("AL18").Value = 90  # was 85
$ws.Range("AM18").Value = 120  # was 100
$ws.Range("AN18").Value = 95  # was 90
$ws.Range("AO18").Value = 6.6  # was 6.4

# Row 20
$ws.Range("F20").Value = 1.39  # was 1.4
$ws.Range("I20").Value = 10  # was 9.4
$ws.Range("J20").Value = 5.3  # was 5.2
$ws.Range("K20").Value = 6  # was 5.8
$ws.Range("T20").Value = 1.92  # was 1.91
$ws.Range("U20").Value = 2  # was 1.95
$ws.Range("W20").Value = 3.35  # was 3.3
$ws.Range("AJ20").Value = 12.5  # was 13.5

# Row 21
$ws.Range("K21").Value = 3.75  # was 3.7
$ws.Range("L21").Value = 1.39  # was 1.35
$ws.Range("M21").Value = 1.06  # was 1.07

# Row 22
$ws.Range("L22").Value = 1.25  # was 1.26
$ws.Range("AD22").Value = 38  # was 36

Write-Host "Applied 234 cell updates."
